# Updates cryptos list price/volume figures (refreshed crypto feed snapshot).
# All "Price" (D) and "Volume(1h)" (E) cells are stored as text in the sheet
# (t="inlineStr"), even when the text looks numeric (e.g. "1.002"). Plain
# `.Value = "1.002"` would make Excel auto-coerce that into a real number, so
# any D-column value that parses as a float is written with a leading `'`
# (the classic "force text" quote-prefix) and the cell style is immediately
# reset to "Normal" so no stray quote-prefix/number-format style lingers on
# the cell - only the text content changes, matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.664.67"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").Value = "1.787.38"
$ws.Range("E3").Value = "  -1.74%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'308.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").Value = "'0.4537"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.94%  "
$ws.Range("D8").Value = "'0.3689"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("D9").Value = "'0.07255"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("D10").Value = "'0.8526"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").Value = "'20.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("D12").Value = "1.790.50"
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("D13").Value = "'6.516"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("D14").Value = "'5.291"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "'0.07029"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").Value = "'90.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.30%  "
$ws.Range("D17").Value = "'1.004"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "'0.000008587"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("D19").Value = "'1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "'14.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.54%  "
$ws.Range("D21").Value = "26.669.42"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").Value = "'5.249"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "'10.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.53%  "
$ws.Range("D24").Value = "2.017.45"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").Value = "'1.906"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.71%  "
$ws.Range("D26").Value = "'150.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("D27").Value = "'2.161"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.24%  "
$ws.Range("D28").Value = "'18.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("D29").Value = "'5.178"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.36%  "
$ws.Range("D30").Value = "'113.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.63%  "
$ws.Range("D31").Value = "'0.08822"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").Value = "'0.7538"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("D33").Value = "'1.151"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("D34").Value = "'4.432"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.82%  "
$ws.Range("D35").Value = "'2.866"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").Value = "'1.003"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("D37").Value = "'1.108"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.65%  "
$ws.Range("D38").Value = "'0.01935"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("D39").Value = "'0.05187"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("D40").Value = "'7.122"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.86%  "
$ws.Range("D41").Value = "'2.857"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.5186"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.01%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'2.312"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.37%  "
$ws.Range("D44").Value = "'0.1643"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.03%  "
$ws.Range("D45").Value = "'8.439"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.94%  "
$ws.Range("D46").Value = "'0.4929"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.98%  "
$ws.Range("D47").Value = "'1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("E48").Value = "  -3.96%  "
$ws.Range("D49").Value = "'103.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("D50").Value = "'1.639"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.92%  "
$ws.Range("D51").Value = "'0.06275"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.36%  "
